$d = $word.ActiveDocument

# Locate the existing "Platform impact" bullet paragraph in the
# "KEY ACHIEVEMENTS AND IMPACT" section - the new bullets are inserted
# right after it (and before the "TECHNICAL SKILLS" heading).
$paras = $d.Paragraphs
$targetIdx = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations*") {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -gt 0) {
    $newBullets = @(
        "• Real-time collaboration at national scale",
        "• Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%",
        "• Increased voter turnout prediction accuracy from 71% to 87%",
        "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
    )

    $prevIdx = $targetIdx
    foreach ($bullet in $newBullets) {
        $prevRange = $d.Paragraphs.Item($prevIdx).Range
        $prevRange.Collapse(0)
        $prevRange.InsertParagraphAfter()
        $prevIdx = $prevIdx + 1
        $newRange = $d.Paragraphs.Item($prevIdx).Range
        $newRange.Collapse(1)
        $newRange.InsertAfter($bullet)
    }
}
